$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, pushing existing data rows down
$ws.Rows.Item(2).Insert()

# Reset formatting on the new row (Insert copies formatting from the row above)
$ws.Rows.Item(2).ClearFormats()

# Populate the new row 2 with data (new weekly price record)
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C2").Value = "Arica y Parinacota"
$ws.Range("D2").Value = 44496
$ws.Range("D2").NumberFormat = $ws.Range("D3").NumberFormat
$ws.Range("E2").Value = 15
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100107
$ws.Range("H2").Value = "Otros"
$ws.Range("I2").Value = 100107002
$ws.Range("J2").Value = "Chirimoya"
$ws.Range("K2").Value = "Cultivar IV Región"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 23000
$ws.Range("O2").Value = 24000
$ws.Range("P2").Value = 23500
$ws.Range("Q2").Value = "$/caja 12 kilos"
$ws.Range("R2").Value = "Región de Coquimbo"
$ws.Range("S2").Value = 1958
$ws.Range("T2").Value = 12
